$wb = $excel.ActiveWorkbook

# --- Insert a new "animals" worksheet between "info" and "addresses" ---
$infoSheet = $wb.Worksheets.Item("info")
$animalsSheet = $wb.Worksheets.Add($null, $infoSheet)
$animalsSheet.Name = "animals"

# Populate the animals sheet. Cells are written in this specific order so
# that new shared-string entries land in the same order as the target file.
$animalsSheet.Range("A7").Value = "Armadillo"
$animalsSheet.Range("A6").Value = "Anteater"
$animalsSheet.Range("A2").Value = "Akita"
$animalsSheet.Range("A3").Value = "Albatross"
$animalsSheet.Range("A1").Value = "animal"
$animalsSheet.Range("A4").Value = "Alpaca"
$animalsSheet.Range("A5").Value = "Anchovy"

# --- Update selection on the "info" sheet (no longer the active tab) ---
$infoSheet.Range("A2").Select()

# --- Make "animals" the active sheet/tab, with row 2 selected ---
$animalsSheet.Activate()
$animalsSheet.Rows(2).Select()

# --- Cosmetic page setup tweak on "info" sheet seen in the target file ---
$infoSheet.PageSetup.Orientation = 1
